# Update radar for S1 2019
# - Apache Kafka & Elasticsearch move from "tools" quadrant to "platforms"
# - ML.NET gets its missing ring ("asses")
# - New entry: JWT (adopt / techniques)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2019-S1")

# Apache Kafka (row 7) and Elasticsearch (row 8): quadrant "tools" -> "platforms"
$ws.Range("C7").Value = "platforms"
$ws.Range("C8").Value = "platforms"

# ML.NET (row 10) was missing its ring value
$ws.Range("B10").Value = "asses"
$ws.Range("B10").WrapText = $true
$ws.Range("B10").VerticalAlignment = -4108

# New row 17: JWT
$ws.Range("A17").Value = "JWT"
$ws.Range("B17").Value = "adopt"
$ws.Range("C17").Value = "techniques"
$ws.Range("D17").Value = $false
$ws.Range("E17").Value = "JSON Web Token is a JSON-based open standard (RFC 7519) for creating access tokens that assert some number of claims."

$ws.Range("A17:D17").WrapText = $true
$ws.Range("A17:D17").VerticalAlignment = -4108
$ws.Range("E17").WrapText = $true
$ws.Rows.Item(17).RowHeight = 29

# Restore view state (best effort)
$ws.Range("D22").Select()
